$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the new backlog item in row 59 (A: task text, B: status "no comenzado")
$ws.Range("A59").Value = "corregir generacion de recibos reporte"
$ws.Range("B59").Value = "no comenzado"

# Update the active view/selection to mirror the scroll down to the new row
$ws.Application.ActiveWindow.ScrollRow = 42
$ws.Range("A60").Select() | Out-Null
